$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F3").Value = -2
$ws.Range("F6").Value = 0
$ws.Range("F10").Value = 5
$ws.Range("F11").Value = 6
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("F28").Value = 3
$ws.Range("F33").Value = 3
$ws.Range("F38").Value = 3
$ws.Range("F42").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("F48").Value = 4
$ws.Range("F49").Value = 3
$ws.Range("F66").Value = -5
$ws.Range("F67").Value = 3
$ws.Range("F70").Value = -6
